$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.551.29'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '2.050.51'
$ws.Range("E3").Value = '  -0.13%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '246.36'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").Value = '  +0.03%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '54.57'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -6.51%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '61.29'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.89%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.365'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.59%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0749'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -4.01%  '
$ws.Range("E12").Value = '  -3.54%  '
$ws.Range("E13").Value = '  +8.22%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '14.70'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.25%  '
$ws.Range("D15").Value = '2.348.60'
$ws.Range("E15").Value = '  -0.22%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '5.44'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -4.22%  '
$ws.Range("D17").Value = '2.059.50'
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '36.460.77'
$ws.Range("E18").Value = '  -1.30%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '17.11'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -5.36%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '71.75'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("E21").Value = '  -4.24%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '237.47'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.21'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -4.14%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("E25").Value = '  -2.76%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +3.56%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '165.20'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.51%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '9.18'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -10.30%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '19.90'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("E30").Value = '  -2.64%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.19'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +7.44%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.05'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -9.06%  '
$ws.Range("E33").Value = '  -5.29%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0592'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.02%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0871'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.34%  '
$ws.Range("B36").Value = 'BinanceUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -1.03%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -4.66%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.04'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.44%  '
$ws.Range("E40").Value = '  -7.20%  '
$ws.Range("E41").Value = '  -5.01%  '
$ws.Range("E42").Value = '  -4.67%  '
$ws.Range("E43").Value = '  -4.75%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '94.12'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.41%  '
$ws.Range("E45").Value = '  -5.15%  '
$ws.Range("D46").Value = '1.404.79'
$ws.Range("E46").Value = '  +7.81%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '15.86'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -6.30%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.41'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +9.33%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.92'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.71%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.26'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.29%  '
$ws.Range("D51").Value = '2.233.81'
$ws.Range("E51").Value = '  -0.28%  '
